$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3
$ws.Range("D2").Value = "Medium"
$ws.Range("E2").Value = "TODO"
$ws.Range("F2").Value = 1

$ws.Range("C3").Value = 5
$ws.Range("D3").Value = "Low"
$ws.Range("E3").Value = "TODO"
$ws.Range("F3").Value = 55

$ws.Range("C4").Value = 5
$ws.Range("D4").Value = "Low"
$ws.Range("E4").Value = "TODO"
$ws.Range("F4").Value = 34

$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "High"
$ws.Range("E5").Value = "TODO"
$ws.Range("F5").Value = 13

$ws.Range("C6").Value = 5
$ws.Range("D6").Value = "Low"
$ws.Range("E6").Value = "TODO"
$ws.Range("F6").Value = 55

$ws.Range("C7").Value = 5
$ws.Range("D7").Value = "Low"
$ws.Range("E7").Value = "TODO"
$ws.Range("F7").Value = 13

$ws.Range("C8").Value = 5
$ws.Range("D8").Value = "Low"
$ws.Range("E8").Value = "TODO"
$ws.Range("F8").Value = 55

$ws.Range("C9").Value = 5
$ws.Range("D9").Value = "Low"
$ws.Range("E9").Value = "TODO"
$ws.Range("F9").Value = 34

$ws.Range("C10").Value = 5
$ws.Range("D10").Value = "Low"
$ws.Range("E10").Value = "TODO"
$ws.Range("F10").Value = 34

$ws.Range("C11").Value = 3
$ws.Range("D11").Value = "High"
$ws.Range("E11").Value = "TODO"
$ws.Range("F11").Value = 13

$ws.Range("C12").Value = 4
$ws.Range("D12").Value = "Medium"
$ws.Range("E12").Value = "TODO"
$ws.Range("F12").Value = 13

$ws.Range("C13").Value = 3
$ws.Range("D13").Value = "High"
$ws.Range("E13").Value = "TODO"
$ws.Range("F13").Value = 21

$ws.Range("C14").Value = 2
$ws.Range("D14").Value = "High"
$ws.Range("E14").Value = "DEVELOPING"
$ws.Range("F14").Value = 34

$ws.Range("C15").Value = 3
$ws.Range("D15").Value = "High"
$ws.Range("E15").Value = "TODO"
$ws.Range("F15").Value = 13

$ws.Range("C16").Value = 3
$ws.Range("D16").Value = "High"
$ws.Range("E16").Value = "TODO"
$ws.Range("F16").Value = 34

$ws.Range("C17").Value = 5
$ws.Range("D17").Value = "Low"
$ws.Range("E17").Value = "TODO"
$ws.Range("F17").Value = 89

$ws.Range("C18").Value = 5
$ws.Range("D18").Value = "Low"
$ws.Range("E18").Value = "TODO"
$ws.Range("F18").Value = 55

$ws.Range("C19").Value = 5
$ws.Range("D19").Value = "Low"
$ws.Range("E19").Value = "TODO"
$ws.Range("F19").Value = 55

$ws.Range("C20").Value = 3
$ws.Range("D20").Value = "Medium"
$ws.Range("E20").Value = "TODO"
$ws.Range("F20").Value = 13

$ws.Range("C21").Value = 2
$ws.Range("D21").Value = "High"
$ws.Range("E21").Value = "TODO"
$ws.Range("F21").Value = 8
